# Updated cryptos list values (Price / Volume(1h) columns) per the commit
# diff. Every D/E data cell in this sheet is stored as plain text (inline
# strings) -- including price values that happen to look numeric, e.g.
# "1.005" or "331.66". Assigning those through .Value as a bare string
# makes Excel silently reinterpret them as real numbers, so numeric-looking
# values below are written with a leading apostrophe to force text.
#
# That apostrophe leaves the cell on a "quote prefix" style, which the
# original cells never had, so afterwards we paste-special (Formats only)
# the plain format from an untouched column-D data cell (D6) onto every
# cell we forced to text, putting them back on the default/unstyled format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.455.63'
$ws.Range('E2').Value = '  -0.35%  '
$ws.Range('D3').Value = '1.826.90'
$ws.Range('E3').Value = '  -1.76%  '
$ws.Range('D4').Value = "'1.005"
$ws.Range('E4').Value = '  -0.57%  '
$ws.Range('D5').Value = "'331.66"
$ws.Range('E5').Value = '  -0.45%  '
$ws.Range('E6').Value = '  -0.72%  '
$ws.Range('D7').Value = "'0.4582"
$ws.Range('E7').Value = '  -1.69%  '
$ws.Range('D8').Value = "'0.3819"
$ws.Range('E8').Value = '  -1.67%  '
$ws.Range('D9').Value = "'46.30"
$ws.Range('E9').Value = '  +1.23%  '
$ws.Range('D10').Value = "'0.07893"
$ws.Range('E10').Value = '  -0.72%  '
$ws.Range('D11').Value = "'0.9696"
$ws.Range('E11').Value = '  -3.06%  '
$ws.Range('D12').Value = "'21.07"
$ws.Range('E12').Value = '  -2.48%  '
$ws.Range('D13').Value = '1.831.62'
$ws.Range('E13').Value = '  -1.61%  '
$ws.Range('D14').Value = "'5.879"
$ws.Range('E14').Value = '  -1.57%  '
$ws.Range('D15').Value = "'7.042"
$ws.Range('E15').Value = '  -2.18%  '
$ws.Range('D16').Value = "'1.005"
$ws.Range('E16').Value = '  -0.86%  '
$ws.Range('D17').Value = "'89.62"
$ws.Range('E17').Value = '  +2.14%  '
$ws.Range('D18').Value = "'0.06630"
$ws.Range('E18').Value = '  -1.18%  '
$ws.Range('E19').Value = '  -1.27%  '
$ws.Range('D20').Value = "'17.12"
$ws.Range('E20').Value = '  +1.33%  '
$ws.Range('E21').Value = '  -0.73%  '
$ws.Range('D22').Value = '27.451.04'
$ws.Range('E22').Value = '  -0.33%  '
$ws.Range('D23').Value = "'5.336"
$ws.Range('E23').Value = '  -1.71%  '
$ws.Range('D24').Value = "'10.81"
$ws.Range('E24').Value = '  -0.14%  '
$ws.Range('D25').Value = "'2.303"
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('D26').Value = '2.045.97'
$ws.Range('E26').Value = '  -1.82%  '
$ws.Range('D27').Value = "'156.65"
$ws.Range('E27').Value = '  -1.33%  '
$ws.Range('D28').Value = "'19.35"
$ws.Range('E28').Value = '  -1.64%  '
$ws.Range('D29').Value = "'2.060"
$ws.Range('E29').Value = '  -2.79%  '
$ws.Range('D30').Value = "'5.275"
$ws.Range('E30').Value = '  -1.75%  '
$ws.Range('E31').Value = '  -2.37%  '
$ws.Range('D32').Value = "'0.9479"
$ws.Range('E32').Value = '  -2.26%  '
$ws.Range('D33').Value = "'0.09325"
$ws.Range('E33').Value = '  -1.30%  '
$ws.Range('D34').Value = "'3.574"
$ws.Range('E34').Value = '  -1.90%  '
$ws.Range('D35').Value = "'5.245"
$ws.Range('E35').Value = '  -0.71%  '
$ws.Range('D36').Value = "'1.327"
$ws.Range('E36').Value = '  +0.13%  '
$ws.Range('D37').Value = "'0.05921"
$ws.Range('E37').Value = '  -1.58%  '
$ws.Range('D38').Value = "'0.02182"
$ws.Range('E38').Value = '  -1.22%  '
$ws.Range('D39').Value = "'1.160"
$ws.Range('E39').Value = '  -2.45%  '
$ws.Range('D40').Value = "'8.021"
$ws.Range('E40').Value = '  -1.85%  '
$ws.Range('D41').Value = "'0.5761"
$ws.Range('E41').Value = '  -2.37%  '
$ws.Range('D42').Value = "'0.1832"
$ws.Range('E42').Value = '  -2.24%  '
$ws.Range('D43').Value = "'10.05"
$ws.Range('E43').Value = '  -1.42%  '
$ws.Range('D44').Value = "'1.264"
$ws.Range('E44').Value = '  +1.91%  '
$ws.Range('D45').Value = "'11.98"
$ws.Range('E45').Value = '  -1.12%  '
$ws.Range('D46').Value = "'0.5447"
$ws.Range('E46').Value = '  -2.84%  '
$ws.Range('D47').Value = "'1.870"
$ws.Range('E47').Value = '  -1.97%  '
$ws.Range('D48').Value = "'0.06610"
$ws.Range('E48').Value = '  -2.11%  '
$ws.Range('D49').Value = "'110.47"
$ws.Range('E49').Value = '  -1.55%  '
$ws.Range('E50').Value = '  -0.75%  '
$ws.Range('D51').Value = "'1.041"
$ws.Range('E51').Value = '  -1.18%  '

# Reset the forced-text cells above back to the default (unstyled) format.
$ws.Range('D6').Copy() | Out-Null
$ws.Range('D4').PasteSpecial(-4122)
$ws.Range('D5').PasteSpecial(-4122)
$ws.Range('D7').PasteSpecial(-4122)
$ws.Range('D8').PasteSpecial(-4122)
$ws.Range('D9').PasteSpecial(-4122)
$ws.Range('D10').PasteSpecial(-4122)
$ws.Range('D11').PasteSpecial(-4122)
$ws.Range('D12').PasteSpecial(-4122)
$ws.Range('D14').PasteSpecial(-4122)
$ws.Range('D15').PasteSpecial(-4122)
$ws.Range('D16').PasteSpecial(-4122)
$ws.Range('D17').PasteSpecial(-4122)
$ws.Range('D18').PasteSpecial(-4122)
$ws.Range('D20').PasteSpecial(-4122)
$ws.Range('D23').PasteSpecial(-4122)
$ws.Range('D24').PasteSpecial(-4122)
$ws.Range('D25').PasteSpecial(-4122)
$ws.Range('D27').PasteSpecial(-4122)
$ws.Range('D28').PasteSpecial(-4122)
$ws.Range('D29').PasteSpecial(-4122)
$ws.Range('D30').PasteSpecial(-4122)
$ws.Range('D32').PasteSpecial(-4122)
$ws.Range('D33').PasteSpecial(-4122)
$ws.Range('D34').PasteSpecial(-4122)
$ws.Range('D35').PasteSpecial(-4122)
$ws.Range('D36').PasteSpecial(-4122)
$ws.Range('D37').PasteSpecial(-4122)
$ws.Range('D38').PasteSpecial(-4122)
$ws.Range('D39').PasteSpecial(-4122)
$ws.Range('D40').PasteSpecial(-4122)
$ws.Range('D41').PasteSpecial(-4122)
$ws.Range('D42').PasteSpecial(-4122)
$ws.Range('D43').PasteSpecial(-4122)
$ws.Range('D44').PasteSpecial(-4122)
$ws.Range('D45').PasteSpecial(-4122)
$ws.Range('D46').PasteSpecial(-4122)
$ws.Range('D47').PasteSpecial(-4122)
$ws.Range('D48').PasteSpecial(-4122)
$ws.Range('D49').PasteSpecial(-4122)
$ws.Range('D51').PasteSpecial(-4122)
$excel.CutCopyMode = 0
